$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Förändrad" (last-changed) date in column C advances by one day (46065 -> 46066)
# for every data row (2-8).
for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r, 3).Value = 46066
}

# Rows 6 and 7 swap their Beteckning (A), Datum (B) and Area ha (G) values.
$a6 = $ws.Range("A6").Value2
$b6 = $ws.Range("B6").Value2
$g6 = $ws.Range("G6").Value2

$a7 = $ws.Range("A7").Value2
$b7 = $ws.Range("B7").Value2
$g7 = $ws.Range("G7").Value2

$ws.Range("A6").Value = $a7
$ws.Range("B6").Value = $b7
$ws.Range("G6").Value = $g7

$ws.Range("A7").Value = $a6
$ws.Range("B7").Value = $b6
$ws.Range("G7").Value = $g6
